$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns are treated as text so that
# numeric-looking strings (e.g. "1.139") are not coerced into numbers,
# matching the original inlineStr cell type.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.169.42'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.868.45'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '311.75'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '0.5049'
$ws.Range("E7").Value = '  -1.62%  '
$ws.Range("D8").Value = '0.3922'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.09706'
$ws.Range("E9").Value = '  -3.30%  '
$ws.Range("D10").Value = '1.139'
$ws.Range("D11").Value = '40.85'
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").Value = '6.508'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").Value = '20.91'
$ws.Range("E13").Value = '  +0.49%  '
$ws.Range("D14").Value = '1.858.85'
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("D16").Value = '7.404'
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").Value = '0.00001129'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = '92.87'
$ws.Range("E18").Value = '  -1.86%  '
$ws.Range("D19").Value = '0.06621'
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("D20").Value = '17.53'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").Value = '6.161'
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("D23").Value = '28.219.92'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '11.36'
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").Value = '2.279'
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("D26").Value = '2.534'
$ws.Range("E26").Value = '  +2.58%  '
$ws.Range("D27").Value = '2.087.29'
$ws.Range("E27").Value = '  +2.40%  '
$ws.Range("D28").Value = '21.24'
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").Value = '158.21'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").Value = '127.41'
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("D31").Value = '0.1061'
$ws.Range("E31").Value = '  -3.08%  '
$ws.Range("D32").Value = '1.066'
$ws.Range("E32").Value = '  -0.48%  '
$ws.Range("D33").Value = '5.630'
$ws.Range("E33").Value = '  -0.35%  '
$ws.Range("D34").Value = '3.622'
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("D35").Value = '9.572'
$ws.Range("E35").Value = '  +4.61%  '
$ws.Range("D36").Value = '0.06719'
$ws.Range("E36").Value = '  -3.03%  '
$ws.Range("D37").Value = '0.02391'
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("D38").Value = '0.2177'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("D39").Value = '11.50'
$ws.Range("E39").Value = '  -1.41%  '
$ws.Range("D40").Value = '0.6354'
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("D41").Value = '4.970'
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("D42").Value = '1.178'
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").Value = '13.57'
$ws.Range("E44").Value = '  +1.46%  '
$ws.Range("D45").Value = '0.6007'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").Value = '3.659'
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("D47").Value = '1.259'
$ws.Range("E47").Value = '  -2.24%  '
$ws.Range("D48").Value = '124.24'
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").Value = '1.195'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").Value = '0.06829'
$ws.Range("E51").Value = '  +0.56%  '

# Reset style to the workbook default so no stray cell-style references
# are introduced (keeps cells as plain/default-styled text cells).
$ws.Range("D2:E51").Style = "Normal"

